# Fix for little mistake in SMILES codes which induced wrong data in
# descriptors and HCA. This updates the recomputed RDKit descriptor values
# for the molecule in row 3 (row index 3, the second data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.507494567271353
$ws.Range("C3").Value = 5.507494567271353
$ws.Range("D3").Value = 0.6440460443436633
$ws.Range("E3").Value = -2.022856355505165
$ws.Range("F3").Value = 0.4690255384985064
$ws.Range("G3").Value = 11.92857142857143
$ws.Range("H3").Value = 274.413
$ws.Range("I3").Value = 255.261
$ws.Range("J3").Value = 274.028479478
$ws.Range("K3").Value = 86
$ws.Range("M3").Value = 0.2468426025619926
$ws.Range("N3").Value = -0.3219679757054398
$ws.Range("O3").Value = 0.3219679757054398
$ws.Range("P3").Value = 0.2468426025619926
$ws.Range("Q3").Value = 1.071428571428571
$ws.Range("R3").Value = 1.642857142857143
$ws.Range("S3").Value = 2.071428571428572
$ws.Range("T3").Value = 32.92115427442042
$ws.Range("U3").Value = 10.90606161098254
$ws.Range("V3").Value = 2.1729804028189
$ws.Range("W3").Value = -2.190680104542174
$ws.Range("X3").Value = 2.603105594315437
$ws.Range("Y3").Value = -1.938755943544446
$ws.Range("Z3").Value = 8.675076046216581
$ws.Range("AA3").Value = 0.2804267459244999
$ws.Range("AB3").Value = 2.134451442561562
$ws.Range("AC3").Value = 3.776797252983614
$ws.Range("AD3").Value = 167.9693269174992
$ws.Range("AE3").Value = 10.86396103067893
$ws.Range("AF3").Value = 9.023988953752012
$ws.Range("AG3").Value = 12.36790588753511
$ws.Range("AH3").Value = 6.681980515339465
$ws.Range("AI3").Value = 4.79499275987393
$ws.Range("AJ3").Value = 10.17852728414391
$ws.Range("AK3").Value = 2.312909769600933
$ws.Range("AL3").Value = 10.13633420366469
$ws.Range("AM3").Value = 1.255622677552617
$ws.Range("AN3").Value = 7.976868250331396
$ws.Range("AO3").Value = 0.7803654957880988
$ws.Range("AP3").Value = 6.642842515024668
$ws.Range("AR3").Value = 986.1165664634418
$ws.Range("AS3").Value = 15.27
$ws.Range("AT3").Value = 9.492720917643227
$ws.Range("AU3").Value = 7.547184073764443
$ws.Range("AV3").Value = 100.8089988065543
$ws.Range("AW3").Value = 9.047494323423635
$ws.Range("BH3").Value = 31.40718409476669
$ws.Range("BJ3").Value = 13.21376392902584
$ws.Range("BK3").Value = 9.047494323423635
$ws.Range("BP3").Value = 20.77121159907187
$ws.Range("BQ3").Value = 30.4723247492662
$ws.Range("BY3").Value = 30.4723247492662
$ws.Range("BZ3").Value = 20.85435041206228
$ws.Range("CB3").Value = 20.77121159907187
$ws.Range("CG3").Value = 18.46
$ws.Range("CI3").Value = 0
$ws.Range("CL3").Value = 13.21376392902584
$ws.Range("CN3").Value = 0
$ws.Range("CO3").Value = 11.38172479611316
$ws.Range("CP3").Value = 25.60935934877237
$ws.Range("CR3").Value = 20.85435041206228
$ws.Range("CS3").Value = 11.01498913454271
$ws.Range("CT3").Value = 8.968239795918366
$ws.Range("CU3").Value = 0
$ws.Range("CW3").Value = -2.022856355505165
$ws.Range("CX3").Value = 3.307730694129504
$ws.Range("DA3").Value = 7.370785619803477
$ws.Range("DD3").Value = 14
$ws.Range("DE3").Value = 0
$ws.Range("DM3").Value = 5
$ws.Range("DN3").Value = 0
$ws.Range("DP3").Value = 9
$ws.Range("DU3").Value = 3.770200000000003
$ws.Range("DV3").Value = 72.99400000000004
